$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Use the sheet's used range to be safe
$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count

for ($r = 2; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is the 7th column
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
